$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.239.28"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "2.237.65"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.28%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "2.579.94"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "2.321.06"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.824"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.02%  "
$ws.Range("D18").Value = "43.974.92"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "0.0₃0960"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.67%  "
$ws.Range("E28").Value = "  +4.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0791"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.76%  "
$ws.Range("E36").Value = "  +1.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.107"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("E38").Value = "  -8.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0297"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.80%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "1.724.78"
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "82.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.190"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.41%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.70%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.04%  "
